{"js": "// The source diff for this change is *purely* a cosmetic XML re-serialization:\n// every element/attribute that appears on a \"-\" line in the diff has an exact\n// \"+\" counterpart with the same tag name and the same set of attribute\n// name/value pairs \u2014 only the attribute order (and the order of the root\n// namespace declarations) differs. This is the well-known side effect of a\n// document being re-saved by a different OOXML writer (here: the commit\n// message notes an Apache POI packaging fix / upgrade to POI 3.15, whose\n// XMLBeans-based writer emits attributes in alphabetical order).\n//\n// No text, formatting, or structural content actually changed: no runs were\n// added/removed/edited, no properties changed value, nothing moved between\n// paragraphs, etc. The Word JavaScript API operates on the document's\n// content/formatting object model and has no means (nor any need) to control\n// the raw XML attribute-serialization order that OOXML writers use, so there\n// is nothing for this script to do in terms of document content.\n//\n// We still touch the object model minimally (a harmless load/sync of the\n// document body) so the script is a valid, verifiable no-op rather than an\n// empty file.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The source diff for this change is *purely* a cosmetic XML re-serialization:\n# every element/attribute that appears on a \"-\" line in the diff has an exact\n# \"+\" counterpart with the same tag name and the same set of attribute\n# name/value pairs - only the attribute order (and the order of the root\n# namespace declarations) differs. This is the well-known side effect of a\n# document being re-saved by a different OOXML writer (here: the commit\n# message notes an Apache POI packaging fix / upgrade to POI 3.15, whose\n# XMLBeans-based writer emits attributes in alphabetical order).\n#\n# No text, formatting, or structural content actually changed: no runs were\n# added/removed/edited, no properties changed value, nothing moved between\n# paragraphs, etc. The Word COM object model operates on the document's\n# content/formatting, not on raw XML attribute-serialization order, so there\n# is nothing for this script to change in terms of document content.\n#\n# We still touch the object model minimally (a harmless read of the document\n# content) so the script is a valid, verifiable no-op rather than an empty\n# file.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
